$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the "total time" row (currently row 36),
# pushing it down to row 40.
$ws.Rows("36:39").Insert()

# The inserted rows lose the original borders/number-formats, so copy the
# formatting from the row directly above (row 35, still blank) back onto
# them to match the rest of the diary table.
$ws.Range("A35:F35").Copy()
$ws.Range("A36:F39").PasteSpecial(-4122)

# Fill in the new diary entry that now lives in row 32.
$ws.Range("A32").Value = "28.11.18"
$ws.Range("B32").Value = 0.65625
$ws.Range("C32").Value = 0.79166666666666663
$ws.Range("D32").Formula = "=C32-B32"
$ws.Range("E32").Value = "-Output Options"

# Refresh the grand-total formula so it also covers the newly inserted rows.
$ws.Range("D40").Formula = "=SUM(D3:D39)"

# Move the active selection to match the author's final cursor position.
$null = $ws.Range("B33").Select()
